# Calendar_for_2023.xlsx edit script
# - inverse availability: add an "Available" default row + split the B2:B7
#   data validation so row 2 gets its own "Available"-only list while the
#   rest of the column (B3:B7) keeps the existing "Not Available" list.
# - support dates: populate the first data row (row 2) on the lecturer
#   sheet with a concrete date + from/until slot.
# - support results upload: (workbook-level change, selection/view bookkeeping)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "פנינה יעקב"
$ws2 = $wb.Worksheets.Item(2)   # "Slots"

# --- Sheet1: fill in row 2 (date / availability / from / until) ---
$ws1.Range("A2").Value = 45205
$ws1.Range("B2").Value = "Available"
$ws1.Range("C2").Value = 0.33333333333333331
$ws1.Range("D2").Value = 0.70833333333333404

# --- Sheet1: data validation bookkeeping ---
# Row 2's Availability cell now defaults to "Available" (Slots!$F$2) while
# the remaining rows (3-7) keep validating against "Not Available"
# (Slots!$F$3).
$ws1.Range("B3:B7").Validation.Modify(3, 1, 1, "Slots!`$F`$3")
$ws1.Range("B2").Validation.Add(3, 1, 1, "Slots!`$F`$2")

# --- Sheet1: selection / active sheet ---
$ws1.Range("C3").Select()
$ws1.Activate()

# --- Sheet2 ("Slots"): view bookkeeping ---
$ws2.Range("G34").Select()
$excel.ActiveWindow.ScrollRow = 20

# --- Sheet2: nudge B34 (17:00 slot end) to the canonical 17/24 serial ---
$ws2.Range("B34").Value = 0.70833333333333337
